$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 158 result columns (bet settled: "Acierto" / profit 0.83)
$ws.Range("G158").Value = "Acierto"
$ws.Range("H158").Value = 0.83

# New rows of data appended by the tracker's automatic update
$data = @(
    @(14870661, "2025-10-14", "Zizou Bergs", "Raphael Collignon", "Gana Zizou Bergs", 1.8),
    @(14880136, "2025-10-15", "Yannick Hanfmann", "Matteo Arnaldi", "Gana Matteo Arnaldi", 1.91),
    @(14881573, "2025-10-14", "Eduardo Ribeiro", "Pedro Boscardin Dias", "Gana Eduardo Ribeiro", 2.38),
    @(14881577, "2025-10-14", "Miguel Tobon", "Murkel Dellien", "Gana Miguel Tobon", 2.63)
)

$startRow = 159
$endRow = $startRow + $data.Count - 1

# Keep the "fecha" column as plain text (matches the rest of the sheet,
# which stores dates as literal strings, not Excel date serials) -
# format as Text BEFORE writing so Excel doesn't auto-convert the value.
$ws.Range("B$startRow`:B$endRow").NumberFormat = "@"

# "resultado" / "profit" start out blank (pending bets, same as every
# other not-yet-settled row already on the sheet).
$ws.Range("G$startRow`:H$endRow").NumberFormat = "@"

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $rec = $data[$i]
    $ws.Cells.Item($row, 1).Value = $rec[0]
    $ws.Cells.Item($row, 2).Value = $rec[1]
    $ws.Cells.Item($row, 3).Value = $rec[2]
    $ws.Cells.Item($row, 4).Value = $rec[3]
    $ws.Cells.Item($row, 5).Value = $rec[4]
    $ws.Cells.Item($row, 6).Value = $rec[5]
}
